# ldlc_suivi_smartphones: append a new price-snapshot column.
#
# Before: ... AL | AM | AN="nom" | AO="url_produit"
# After:  ... AL | AM | AN=<new timestamp col> | AO="nom" | AP="url_produit"
#
# The new "AN" column carries forward the most recent known price (the
# value that was in column AM) for every product row that had a price,
# and stays blank for rows that had no price in AM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at AN, pushing the existing "nom" / "url_produit"
# columns one place to the right (AN->AO, AO->AP). Excel.Range.Insert on an
# entire-column range defaults to shifting the existing columns right.
$ws.Columns("AN:AN").Insert()

# Header for the newly inserted column: the new snapshot timestamp.
$ws.Cells.Item(1, 40).Value2 = "2026-01-29 11:23:22"

# Determine the last used row from the sheet's dimensions.
$lastRow = $ws.UsedRange.Rows.Count

# For every data row, carry the last price (column AM = column 39) forward
# into the freshly-inserted column AN (column 40) when that price exists.
for ($r = 2; $r -le $lastRow; $r++) {
    $lastPrice = $ws.Cells.Item($r, 39).Value2
    if ($lastPrice -ne "") {
        $ws.Cells.Item($r, 40).Value2 = $lastPrice
    }
}
